# Fixed some bugs in CalcAdjacentPay
# Re-orders the per-symbol pay rows (A2:F25) on the active worksheet to the
# corrected order/values produced after the fix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(401, 9, 48, 67, 75, 45),
    @(701, 3, 90, 45, 97, 15),
    @(1202, 2, 10, 10, 10, 10),
    @(101, 9, 30, 15, 60, 15),
    @(1001, 18, 30, 75, 60, 72),
    @(301, 6, 45, 30, 60, 45),
    @(501, 9, 52, 30, 75, 45),
    @(801, 3, 67, 65, 52, 45),
    @(1201, 2, 10, 10, 10, 10),
    @(1203, 3, 15, 15, 15, 15),
    @(902, 1, 0, 0, 0, 0),
    @(201, 9, 30, 15, 45, 30),
    @(601, 9, 60, 67, 60, 42),
    @(901, 16, 15, 45, 60, 60),
    @(1, 0, 2, 2, 2, 2),
    @(502, 0, 4, 0, 0, 0),
    @(802, 0, 4, 5, 4, 0),
    @(2, 0, 2, 2, 2, 2),
    @(3, 0, 3, 3, 3, 3),
    @(1101, 0, 15, 30, 30, 0),
    @(402, 0, 0, 4, 0, 0),
    @(602, 0, 0, 4, 0, 9),
    @(702, 0, 0, 0, 4, 0),
    @(1002, 0, 0, 0, 0, 9)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Count; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $values[$j]
    }
}
